# Weekly fruit/vegetable update: insert a new weekly record at row 210
# (pushing the existing rows 210:279 down to 211:280) and populate the
# new row with this week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 210; Excel shifts rows 210:279 -> 211:280
$ws.Rows.Item(210).Insert()

# Populate the newly inserted row 210 with the new weekly record.
$ws.Cells.Item(210, 1).Value = 3
$ws.Cells.Item(210, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(210, 3).Value = "Coquimbo"
$ws.Cells.Item(210, 4).Value = 44627
$ws.Cells.Item(210, 5).Value = 5
$ws.Cells.Item(210, 6).Value = 100112039
$ws.Cells.Item(210, 7).Value = "Ciboulette"
$ws.Cells.Item(210, 8).Value = "Sin especificar"
$ws.Cells.Item(210, 9).Value = "Primera"
$ws.Cells.Item(210, 10).Value = 45
$ws.Cells.Item(210, 11).Value = 2000
$ws.Cells.Item(210, 12).Value = 2000
$ws.Cells.Item(210, 13).Value = 2000
$ws.Cells.Item(210, 14).Value = "$/docena de atados"
$ws.Cells.Item(210, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(210, 16).Value = 667
$ws.Cells.Item(210, 17).Value = 3
$ws.Cells.Item(210, 18).Value = "Hortaliza"
